# Update example input files to include Dp
#
# On the "ions" worksheet, insert two new columns ("Dp" and "Dp_units")
# immediately before the existing "conc_units" column (which shifts from
# column I to column K).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("ions")

# Insert two blank columns at I:J - this pushes the former column I
# ("conc_units" + its values) to column K.
$ws.Columns("I:J").Insert()

# --- Header row ---
$ws.Range("I1").Value = "Dp"
$ws.Range("J1").Value = "Dp_units"

# --- Data rows: Dp value (column I) and Dp_units (column J) ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "cm^2/s"

$ws.Range("I3").Value = 0.000002
$ws.Range("I3").NumberFormat = "0.00E+00"
$ws.Range("J3").Value = "cm^2/s"

$ws.Range("I4").Value = 0.000002
$ws.Range("I4").NumberFormat = "0.00E+00"
$ws.Range("J4").Value = "cm^2/s"

$ws.Range("I5").Value = 0.000002
$ws.Range("I5").NumberFormat = "0.00E+00"
$ws.Range("J5").Value = "cm^2/s"

$ws.Range("I6").Value = 0.000002
$ws.Range("I6").NumberFormat = "0.00E+00"
$ws.Range("J6").Value = "cm^2/s"
